$d = $word.ActiveDocument

$d.Content.Find.Execute(
    "keep moving with the same but opposite", $true, $false, $false, $false, $false,
    $true, 1, $false,
    "changent la direction de bouger, mais avec la même", 2)

$d.Content.Find.Execute(
    "so again what are the precise positions", $true, $false, $false, $false, $false,
    $true, 1, $false,
    "Donc quelles sont les positions précises", 2)

$d.Content.Find.Execute(
    "where I should place the two ants in", $true, $false, $false, $false, $false,
    $true, 1, $false,
    "où je devrais placer les deux fourmis", 2)

$d.Content.Find.Execute(
    "order to get the longest time before the", $true, $false, $false, $false, $false,
    $true, 1, $false,
    "pour obtenir le temps le plus long avant que", 2)

$d.Content.Find.Execute(
    "last ant falls? The second puzzle is", $true, $false, $false, $false, $false,
    $true, 1, $false,
    "la dernière fourmis tombe? The second puzzle is", 2)
